$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing row (83) for columns C, D, E, F, G
# onto the new row (84) so that wrap-text / font styling matches the rest of
# the table, then set the actual cell values for the new "oil_imports" entry.
$ws.Range("C83:G83").Copy()
$ws.Range("C84").PasteSpecial(-4122)

$ws.Range("A84").Value = "oil_imports"
$ws.Range("B84").Value = "csv"
$ws.Range("C84").Value = "Socioeconomic Data"
$ws.Range("D84").Value = "monthly oil and gas imports"
$ws.Range("E84").Value = "na"
$ws.Range("F84").Value = "Indonesia"
$ws.Range("G84").Value = "https://drive.google.com/drive/folders/1gg1jPYMPD0pWS5mMMmnUMTnYY39pWnFI"

$excel.CutCopyMode = $false

# Update the window view to mirror the author's scroll/selection position.
[void]$ws.Range("G92").Select()
$excel.ActiveWindow.ScrollRow = 80
